# Updates the cryptos price/volume table with refreshed values (GitHub Actions style refresh).
# Cells whose new text would otherwise be auto-parsed by Excel as a number (e.g. "219.39")
# are written with a leading apostrophe so they stay plain text, matching the source data
# which always stores these columns as text.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.550.11'
$ws.Range('E2').Value = '  -2.80%  '
$ws.Range('D3').Value = '1.670.16'
$ws.Range('E3').Value = '  -2.26%  '
$ws.Range('E4').Value = '  +0.38%  '
$ws.Range('D5').Value = '''219.39'
$ws.Range('E5').Value = '  -1.85%  '
$ws.Range('D6').Value = '''0.5132'
$ws.Range('E6').Value = '  -3.34%  '
$ws.Range('E7').Value = '  +0.30%  '
$ws.Range('D8').Value = '''0.06450'
$ws.Range('D9').Value = '''0.2564'
$ws.Range('E9').Value = '  -3.46%  '
$ws.Range('D10').Value = '''19.92'
$ws.Range('E10').Value = '  -4.65%  '
$ws.Range('D11').Value = '''0.07657'
$ws.Range('E11').Value = '  +0.04%  '
$ws.Range('B12').Value = 'WrappedEther'
$ws.Range('C12').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D12').Value = '1.688.86'
$ws.Range('E12').Value = '  -1.59%  '
$ws.Range('B13').Value = 'Polkadot'
$ws.Range('C13').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D13').Value = '''4.339'
$ws.Range('E13').Value = '  -5.39%  '
$ws.Range('B14').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C14').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D14').Value = '1.898.93'
$ws.Range('E14').Value = '  -2.30%  '
$ws.Range('D15').Value = '''0.5559'
$ws.Range('E15').Value = '  -3.20%  '
$ws.Range('D16').Value = '0.0₅8015'
$ws.Range('E16').Value = '  -2.15%  '
$ws.Range('D17').Value = '''64.66'
$ws.Range('E17').Value = '  -4.36%  '
$ws.Range('D18').Value = '26.568.20'
$ws.Range('E18').Value = '  -2.71%  '
$ws.Range('E19').Value = '  +0.33%  '
$ws.Range('D20').Value = '''210.21'
$ws.Range('E20').Value = '  -2.95%  '
$ws.Range('D21').Value = '''4.438'
$ws.Range('E21').Value = '  -5.13%  '
$ws.Range('E22').Value = '  -3.29%  '
$ws.Range('D23').Value = '''5.886'
$ws.Range('E23').Value = '  -1.61%  '
$ws.Range('D24').Value = '''1.007'
$ws.Range('E24').Value = '  +0.41%  '
$ws.Range('D25').Value = '''142.71'
$ws.Range('E25').Value = '  +0.18%  '
$ws.Range('D26').Value = '''1.719'
$ws.Range('E26').Value = '  -1.61%  '
$ws.Range('D27').Value = '''0.1165'
$ws.Range('E27').Value = '  -4.24%  '
$ws.Range('D28').Value = '''6.978'
$ws.Range('E28').Value = '  -4.01%  '
$ws.Range('E29').Value = '  -4.14%  '
$ws.Range('D30').Value = '''0.05198'
$ws.Range('E30').Value = '  -3.35%  '
$ws.Range('E31').Value = '  -2.29%  '
$ws.Range('E32').Value = '  -4.60%  '
$ws.Range('D33').Value = '''3.199'
$ws.Range('E33').Value = '  -6.59%  '
$ws.Range('E34').Value = '  -3.90%  '
$ws.Range('D35').Value = '''2.756'
$ws.Range('E35').Value = '  -4.29%  '
$ws.Range('D36').Value = '''2.375'
$ws.Range('E36').Value = '  -1.97%  '
$ws.Range('D37').Value = '''0.9227'
$ws.Range('E37').Value = '  -2.75%  '
$ws.Range('D38').Value = '''0.5709'
$ws.Range('E38').Value = '  -2.57%  '
$ws.Range('D39').Value = '1.149.99'
$ws.Range('E39').Value = '  +10.38%  '
$ws.Range('D40').Value = '''0.01587'
$ws.Range('E40').Value = '  -2.88%  '
$ws.Range('E41').Value = '  +0.32%  '
$ws.Range('D42').Value = '''0.8315'
$ws.Range('E42').Value = '  -1.08%  '
$ws.Range('D43').Value = '''5.644'
$ws.Range('D44').Value = '''99.88'
$ws.Range('E44').Value = '  -1.19%  '
$ws.Range('D45').Value = '1.808.84'
$ws.Range('E45').Value = '  -2.24%  '
$ws.Range('D46').Value = '0.0₈112'
$ws.Range('E46').Value = '  -1.77%  '
$ws.Range('D47').Value = '''0.4493'
$ws.Range('E47').Value = '  -0.22%  '
$ws.Range('D48').Value = '''55.53'
$ws.Range('D49').Value = '''1.005'
$ws.Range('E49').Value = '  -0.30%  '
$ws.Range('D50').Value = '''7.907'
$ws.Range('E50').Value = '  -2.18%  '
$ws.Range('D51').Value = '''0.05137'
$ws.Range('E51').Value = '  -2.01%  '
